$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 3.184887070376469
$ws.Range("D2").Value = 8.907179426559555
$ws.Range("E2").Value = 13.43293486100962
$ws.Range("F2").Value = 34.97279954082565
$ws.Range("G2").Value = 3.626802777701949
$ws.Range("I2").Value = 24.17107080665094
$ws.Range("J2").Value = 9.875322075795385
$ws.Range("O2").Value = 26.24867479804472
$ws.Range("C3").Value = 3.173039412892372
$ws.Range("D3").Value = 8.918618001936702
$ws.Range("E3").Value = 13.41202101459874
$ws.Range("F3").Value = 34.56710278368461
$ws.Range("G3").Value = 3.630387562387835
$ws.Range("I3").Value = 23.91486072157212
$ws.Range("J3").Value = 9.878801804244505
$ws.Range("O3").Value = 25.97430490216361
$ws.Range("C4").Value = 3.165543185052686
$ws.Range("D4").Value = 8.926895744457138
$ws.Range("E4").Value = 13.40211214639553
$ws.Range("F4").Value = 34.32608609072018
$ws.Range("G4").Value = 3.632704804535959
$ws.Range("I4").Value = 23.76328479053291
$ws.Range("J4").Value = 9.88276646323262
$ws.Range("O4").Value = 25.8120429835944
$ws.Range("C5").Value = 3.162432947957929
$ws.Range("D5").Value = 8.930584414495634
$ws.Range("E5").Value = 13.39881351554922
$ws.Range("F5").Value = 34.23001325366262
$ws.Range("G5").Value = 3.633678420639572
$ws.Range("I5").Value = 23.70302570273661
$ws.Range("J5").Value = 9.884841610473202
$ws.Range("O5").Value = 25.74754968645196
$ws.Range("C6").Value = 3.161913154073516
$ws.Range("D6").Value = 8.931215965042892
$ws.Range("E6").Value = 13.39831047878087
$ws.Range("F6").Value = 34.21419303463604
$ws.Range("G6").Value = 3.633841862984083
$ws.Range("I6").Value = 23.69311276154824
$ws.Range("J6").Value = 9.885213932416457
$ws.Range("O6").Value = 25.73694103359404
$ws.Range("C7").Value = 3.165501463507393
$ws.Range("D7").Value = 8.926944214058942
$ws.Range("E7").Value = 13.40206466430567
$ws.Range("F7").Value = 34.32478159746022
$ws.Range("G7").Value = 3.632717816198798
$ws.Range("I7").Value = 23.76246591904307
$ws.Range("J7").Value = 9.882792589188515
$ws.Range("O7").Value = 25.81116651585117
$ws.Range("C8").Value = 3.180848187150833
$ws.Range("D8").Value = 8.910863099260368
$ws.Range("E8").Value = 13.42511651026994
$ws.Range("F8").Value = 34.83130573294215
$ws.Range("G8").Value = 3.628014766264922
$ws.Range("I8").Value = 24.08157950035088
$ws.Range("J8").Value = 9.87614240908356
$ws.Range("O8").Value = 26.15282793226649
$ws.Range("C9").Value = 3.209171912478049
$ws.Range("D9").Value = 8.889280952871653
$ws.Range("E9").Value = 13.49347747762701
$ws.Range("F9").Value = 35.88364097495468
$ws.Range("G9").Value = 3.619708880248593
$ws.Range("I9").Value = 24.74981802204094
$ws.Range("J9").Value = 9.877608563381783
$ws.Range("O9").Value = 26.8687921675616
$ws.Range("C10").Value = 3.228886812376733
$ws.Range("D10").Value = 8.879490553274922
$ws.Range("E10").Value = 13.55763621742071
$ws.Range("F10").Value = 36.68574067888086
$ws.Range("G10").Value = 3.614158490422723
$ws.Range("I10").Value = 25.26237089453181
$ws.Range("J10").Value = 9.887524230839666
$ws.Range("O10").Value = 27.41830099875531
$ws.Range("C11").Value = 3.237614892092953
$ws.Range("D11").Value = 8.876352657301947
$ws.Range("E11").Value = 13.58979817134243
$ws.Range("F11").Value = 37.05538728041412
$ws.Range("G11").Value = 3.611751841278708
$ws.Range("I11").Value = 25.49928966524189
$ws.Range("J11").Value = 9.893949522941307
$ws.Range("O11").Value = 27.67238650064728
$ws.Range("C12").Value = 3.240885143120051
$ws.Range("D12").Value = 8.8753534292153
$ws.Range("E12").Value = 13.60239967141665
$ws.Range("F12").Value = 37.19591936505108
$ws.Range("G12").Value = 3.610857396270092
$ws.Range("I12").Value = 25.58946456604475
$ws.Range("J12").Value = 9.896657164325662
$ws.Range("O12").Value = 27.76910816928887
$ws.Range("C13").Value = 3.240182392637692
$ws.Range("D13").Value = 8.875560227432794
$ws.Range("E13").Value = 13.5996670164717
$ws.Range("F13").Value = 37.16563065336729
$ws.Range("G13").Value = 3.611049280982573
$ws.Range("I13").Value = 25.57002466897509
$ws.Range("J13").Value = 9.896061830545133
$ws.Range("O13").Value = 27.74825633861128
$ws.Range("C14").Value = 3.237884638579152
$ws.Range("D14").Value = 8.876266663110588
$ws.Range("E14").Value = 13.59082647073775
$ws.Range("F14").Value = 37.06693856223235
$ws.Range("G14").Value = 3.61167791662476
$ws.Range("I14").Value = 25.5066996760724
$ws.Range("J14").Value = 9.894166785621948
$ws.Range("O14").Value = 27.68033422296784
$ws.Range("C15").Value = 3.23647264543239
$ws.Range("D15").Value = 8.876723985602377
$ws.Range("E15").Value = 13.58546622723222
$ws.Range("F15").Value = 37.00655521686863
$ws.Range("G15").Value = 3.612065172022336
$ws.Range("I15").Value = 25.46796864448457
$ws.Range("J15").Value = 9.893041739964753
$ws.Range("O15").Value = 27.63879317500574
$ws.Range("C16").Value = 3.228311548051888
$ws.Range("D16").Value = 8.879722084771453
$ws.Range("E16").Value = 13.55559378609663
$ws.Range("F16").Value = 36.66166815196527
$ws.Range("G16").Value = 3.614318141252795
$ws.Range("I16").Value = 25.24695633307855
$ws.Range("J16").Value = 9.887142786062588
$ws.Range("O16").Value = 27.40177127560219
$ws.Range("C17").Value = 3.223243194034743
$ws.Range("D17").Value = 8.881898214503742
$ws.Range("E17").Value = 13.53802649215893
$ws.Range("F17").Value = 36.45121898641334
$ws.Range("G17").Value = 3.61573047797931
$ws.Range("I17").Value = 25.11227666083741
$ws.Range("J17").Value = 9.884013840549324
$ws.Range("O17").Value = 27.25735771348779
$ws.Range("C18").Value = 3.22030540942308
$ws.Range("D18").Value = 8.883273731021669
$ws.Range("E18").Value = 13.52820271101449
$ws.Range("F18").Value = 36.33063071707133
$ws.Range("G18").Value = 3.616553953218653
$ws.Range("I18").Value = 25.03517087522831
$ws.Range("J18").Value = 9.882394425402193
$ws.Range("O18").Value = 27.17468699768396
$ws.Range("C19").Value = 3.219306854133521
$ws.Range("D19").Value = 8.883760735919116
$ws.Range("E19").Value = 13.52492487384943
$ws.Range("F19").Value = 36.28988400256065
$ws.Range("G19").Value = 3.616834683653743
$ws.Range("I19").Value = 25.00912821015884
$ws.Range("J19").Value = 9.881877102131957
$ws.Range("O19").Value = 27.14676608934437
$ws.Range("C20").Value = 3.22378507050593
$ws.Range("D20").Value = 8.881653743502332
$ws.Range("E20").Value = 13.5398675715803
$ws.Range("F20").Value = 36.47357536012329
$ws.Range("G20").Value = 3.615578980493172
$ws.Range("I20").Value = 25.12657703735449
$ws.Range("J20").Value = 9.884328269894416
$ws.Range("O20").Value = 27.27269082254469
$ws.Range("C21").Value = 3.238560493944145
$ws.Range("D21").Value = 8.876054037480355
$ws.Range("E21").Value = 13.59341173398039
$ws.Range("F21").Value = 37.095912796592
$ws.Range("G21").Value = 3.61149281323827
$ws.Range("I21").Value = 25.52528795337733
$ws.Range("J21").Value = 9.894715963075924
$ws.Range("O21").Value = 27.70027157568704
$ws.Range("C22").Value = 3.248013611457706
$ws.Range("D22").Value = 8.87349596585937
$ws.Range("E22").Value = 13.63086546988134
$ws.Range("F22").Value = 37.50582122806088
$ws.Range("G22").Value = 3.608920729752565
$ws.Range("I22").Value = 25.78850597668636
$ws.Range("J22").Value = 9.903104544899231
$ws.Range("O22").Value = 27.98262404594374
$ws.Range("C23").Value = 3.242987029473825
$ws.Range("D23").Value = 8.874760529569647
$ws.Range("E23").Value = 13.61065260601587
$ws.Range("F23").Value = 37.2867977169572
$ws.Range("G23").Value = 3.610284523830448
$ws.Range("I23").Value = 25.64780700570682
$ws.Range("J23").Value = 9.898481347597158
$ws.Range("O23").Value = 27.83168995744551
$ws.Range("C24").Value = 3.223540162627329
$ws.Range("D24").Value = 8.881763881248714
$ws.Range("E24").Value = 13.53903435982338
$ws.Range("F24").Value = 36.46346678210745
$ws.Range("G24").Value = 3.615647436664228
$ws.Range("I24").Value = 25.12011082369087
$ws.Range("J24").Value = 9.884185557348077
$ws.Range("O24").Value = 27.26575761288615
$ws.Range("C25").Value = 3.201700739316453
$ws.Range("D25").Value = 8.894053946642025
$ws.Range("E25").Value = 13.47252147058632
$ws.Range("F25").Value = 35.59335385997625
$ws.Range("G25").Value = 3.62185841508966
$ws.Range("I25").Value = 24.56493140925225
$ws.Range("J25").Value = 9.875658675796123
$ws.Range("O25").Value = 26.67064358845211
